{"js": "// Apply the four Russian-translation text edits to GRAPHICS.docx.\n\nconst replacements = [\n  {\n    find: \"\u0420\u0435\u043a\u043e\u043c\u0435\u043d\u0434\u0443\u0435\u043c\u044b\u0435 \u043f\u0430\u0440\u0430\u043c\u0435\u0442\u0440\u044b \u0434\u043b\u044f \u0433\u0440\u0430\u0444\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u043e\u0432\",\n    replace: \"\u0420\u0435\u043a\u043e\u043c\u0435\u043d\u0434\u0443\u0435\u043c\u044b\u0435 \u043f\u0430\u0440\u0430\u043c\u0435\u0442\u0440\u044b \u0433\u0440\u0430\u0444\u0438\u043a\u0438\",\n  },\n  {\n    find:\n      \"Standardization of colours will go a long way to enforce a reliable and positive impression to our identity in the blockchain space.\",\n    replace:\n      \"\u0421\u0442\u0430\u043d\u0434\u0430\u0440\u0442\u0438\u0437\u0430\u0446\u0438\u044f \u0446\u0432\u0435\u0442\u043e\u0432 \u0438\u0433\u0440\u0430\u0435\u0442 \u0432\u0430\u0436\u043d\u0443\u044e \u0440\u043e\u043b\u044c \u0432 \u043f\u043e\u0434\u0434\u0435\u0440\u0436\u0430\u043d\u0438\u0438 \u043d\u0430\u0448\u0435\u0439 \u0438\u0434\u0435\u043d\u0442\u0438\u0447\u043d\u043e\u0441\u0442\u0438 \u0432 \u043a\u0440\u0438\u043f\u0442\u043e\u0432\u0430\u043b\u044e\u0442\u043d\u043e\u043c \u043f\u0440\u043e\u0441\u0442\u0440\u0430\u043d\u0441\u0442\u0432\u0435.\",\n  },\n  {\n    find: \"Official font is\\u00A0\",\n    replace: \"\u041e\u0444\u0438\u0446\u0438\u0430\u043b\u044c\u043d\u044b\u0439 \u0448\u0440\u0438\u0444\u0442\",\n  },\n  {\n    find: \"\\u00A0/ PT Sans\",\n    replace: \"/ PT Sans\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$nbsp = [char]0x00A0\n\n$replacements = @(\n    @{ Find = \"\u0420\u0435\u043a\u043e\u043c\u0435\u043d\u0434\u0443\u0435\u043c\u044b\u0435 \u043f\u0430\u0440\u0430\u043c\u0435\u0442\u0440\u044b \u0434\u043b\u044f \u0433\u0440\u0430\u0444\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u043e\u0432\"; Replace = \"\u0420\u0435\u043a\u043e\u043c\u0435\u043d\u0434\u0443\u0435\u043c\u044b\u0435 \u043f\u0430\u0440\u0430\u043c\u0435\u0442\u0440\u044b \u0433\u0440\u0430\u0444\u0438\u043a\u0438\" },\n    @{ Find = \"Standardization of colours will go a long way to enforce a reliable and positive impression to our identity in the blockchain space.\"; Replace = \"\u0421\u0442\u0430\u043d\u0434\u0430\u0440\u0442\u0438\u0437\u0430\u0446\u0438\u044f \u0446\u0432\u0435\u0442\u043e\u0432 \u0438\u0433\u0440\u0430\u0435\u0442 \u0432\u0430\u0436\u043d\u0443\u044e \u0440\u043e\u043b\u044c \u0432 \u043f\u043e\u0434\u0434\u0435\u0440\u0436\u0430\u043d\u0438\u0438 \u043d\u0430\u0448\u0435\u0439 \u0438\u0434\u0435\u043d\u0442\u0438\u0447\u043d\u043e\u0441\u0442\u0438 \u0432 \u043a\u0440\u0438\u043f\u0442\u043e\u0432\u0430\u043b\u044e\u0442\u043d\u043e\u043c \u043f\u0440\u043e\u0441\u0442\u0440\u0430\u043d\u0441\u0442\u0432\u0435.\" },\n    @{ Find = \"Official font is$nbsp\"; Replace = \"\u041e\u0444\u0438\u0446\u0438\u0430\u043b\u044c\u043d\u044b\u0439 \u0448\u0440\u0438\u0444\u0442\" },\n    @{ Find = \"$nbsp/ PT Sans\"; Replace = \"/ PT Sans\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = 1\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
